# "Generate Report for Handback" — the handback XLIFF files have come back
# from the localization vendor for zh-cn and de-de, so this run:
#   1. Flips the Overview status column for both languages from
#      "Ready for handoff" to "Handed back: in sync with en-US".
#   2. Fills in the (previously empty) Latest Target File / Latest Handback
#      File columns on the per-language sheets, with the target-file cell
#      turned into a hyperlink (mirroring the existing Source File Name link).
#   3. Updates the Latest Handback DateTime for de-de to the new timestamp
#      (zh-cn keeps referring to the already-present handback datetime).
#   4. Widens the Status / Latest Target File / Latest Handback File columns
#      so the longer text fits.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"
$mdName = "7d0caaa4-784d-4383-8fcd-afeb58712726.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3f39f3675436a09efb25e6dd60eedb42a50982f7/e2e/7d0caaa4-784d-4383-8fcd-afeb58712726.md"

# --- 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
# All four cells (Overview!E2/F2 status-per-language, and the Status column on
# each language sheet) shared the same text, so update all of them together.
$overview.Range("E2").Value2 = $newStatus
$overview.Range("F2").Value2 = $newStatus
$zhcn.Range("C2").Value2 = $newStatus
$dede.Range("C2").Value2 = $newStatus

# --- 2. Latest Target File (I2) + Latest Handback File (J2) on zh-cn ---
$zhcn.Range("I2").Value2 = $mdName
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdUrl, "", "", $mdName) | Out-Null
$zhcn.Range("J2").Value2 = "7d0caaa4-784d-4383-8fcd-afeb58712726.a8b938be2a9090e4993f6e7ff7df2b3297997b78.zh-cn.xlf"

# --- 2b. Latest Target File (I2) + Latest Handback File (J2) on de-de ---
$dede.Range("I2").Value2 = $mdName
$dede.Hyperlinks.Add($dede.Range("I2"), $mdUrl, "", "", $mdName) | Out-Null
$dede.Range("J2").Value2 = "7d0caaa4-784d-4383-8fcd-afeb58712726.a8b938be2a9090e4993f6e7ff7df2b3297997b78.de-de.xlf"

# --- 3. Latest Handback DateTime on de-de updates to the new handback run ---
$dede.Range("K2").Value2 = "2016-09-06 03:06:42"

# zh-cn's Latest Handback DateTime (K2) keeps pointing at the same stamp, but
# that stamp itself moves forward to reflect this handback run.
$zhcn.Range("K2").Value2 = "2016-09-06 03:06:35"

# --- 4. Column widths: widen Status + Target/Handback File columns ---
$overview.Columns.Item(5).ColumnWidth = 29.16666666666667   # E: zh-cn status
$overview.Columns.Item(6).ColumnWidth = 29.16666666666667   # F: de-de status

foreach ($ws in @($zhcn, $dede)) {
    $ws.Columns.Item(3).ColumnWidth = 29.16666666666667   # C: Status
    $ws.Columns.Item(9).ColumnWidth = 39.16666666666667   # I: Latest Target File
    $ws.Columns.Item(10).ColumnWidth = 39.16666666666667  # J: Latest Handback File
}
